$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NextBus3")

# Row 2
$ws.Cells.Item(2, 1).Value = "NextBus3"
$ws.Cells.Item(2, 2).Value = 52
$ws.Cells.Item(2, 3).Value = 53009
$ws.Cells.Item(2, 4).Value = "Bishan Int"
$ws.Cells.Item(2, 5).Value = "SBST"
$ws.Cells.Item(2, 6).Value = 45684.40243055556
$ws.Cells.Item(2, 7).Value = 53009
$ws.Cells.Item(2, 8).Value = "WAB"
$ws.Cells.Item(2, 9).Value = "SEA"
$ws.Cells.Item(2, 10).Value = 1
$ws.Cells.Item(2, 11).Value = 28009
$ws.Cells.Item(2, 12).Value = "SD"
$ws.Cells.Item(2, 13).Value = 12101
$ws.Cells.Item(2, 14).Value = 612
$ws.Cells.Item(2, 15).Value = 15
$ws.Cells.Item(2, 16).Value = 613
$ws.Cells.Item(2, 17).Value = 16
$ws.Cells.Item(2, 18).Value = 619
$ws.Cells.Item(2, 19).Value = 15
$ws.Cells.Item(2, 20).Value = 12101
$ws.Cells.Item(2, 21).Value = "Ngee Ann Poly"

# Row 3
$ws.Cells.Item(3, 1).Value = "NextBus3"
$ws.Cells.Item(3, 2).Value = 184
$ws.Cells.Item(3, 3).Value = 44989
$ws.Cells.Item(3, 4).Value = "Gali Batu Ter"
$ws.Cells.Item(3, 5).Value = "SMRT"
$ws.Cells.Item(3, 6).Value = 45684.39763888889
$ws.Cells.Item(3, 7).Value = 44989
$ws.Cells.Item(3, 8).Value = "WAB"
$ws.Cells.Item(3, 9).Value = "SEA"
$ws.Cells.Item(3, 10).Value = 1
$ws.Cells.Item(3, 11).Value = 44989
$ws.Cells.Item(3, 12).Value = "DD"
$ws.Cells.Item(3, 13).Value = 12109
$ws.Cells.Item(3, 14).Value = 622
$ws.Cells.Item(3, 15).Value = 2331
$ws.Cells.Item(3, 16).Value = 621
$ws.Cells.Item(3, 17).Value = 2331
$ws.Cells.Item(3, 18).Value = 602
$ws.Cells.Item(3, 19).Value = 2327
$ws.Cells.Item(3, 20).Value = 12109
$ws.Cells.Item(3, 21).Value = "Opp Ngee Ann Poly"

# Row 4
$ws.Cells.Item(4, 1).Value = "NextBus3"
$ws.Cells.Item(4, 2).Value = 75
$ws.Cells.Item(4, 3).Value = 44989
$ws.Cells.Item(4, 4).Value = "Gali Batu Ter"
$ws.Cells.Item(4, 5).Value = "SMRT"
$ws.Cells.Item(4, 6).Value = 45684.39942129629
$ws.Cells.Item(4, 7).Value = 44989
$ws.Cells.Item(4, 8).Value = "WAB"
$ws.Cells.Item(4, 9).Value = "SEA"
$ws.Cells.Item(4, 10).Value = 1
$ws.Cells.Item(4, 11).Value = 10009
$ws.Cells.Item(4, 12).Value = "SD"
$ws.Cells.Item(4, 13).Value = 12101
$ws.Cells.Item(4, 14).Value = 640
$ws.Cells.Item(4, 15).Value = 32
$ws.Cells.Item(4, 16).Value = 704
$ws.Cells.Item(4, 17).Value = 30
$ws.Cells.Item(4, 18).Value = 633
$ws.Cells.Item(4, 19).Value = 31
$ws.Cells.Item(4, 20).Value = 12101
$ws.Cells.Item(4, 21).Value = "Ngee Ann Poly"

# Row 5
$ws.Cells.Item(5, 1).Value = "NextBus3"
$ws.Cells.Item(5, 2).Value = 184
$ws.Cells.Item(5, 3).Value = 44989
$ws.Cells.Item(5, 4).Value = "Gali Batu Ter"
$ws.Cells.Item(5, 5).Value = "SMRT"
$ws.Cells.Item(5, 6).Value = 45684.40603009259
$ws.Cells.Item(5, 7).Value = 44989
$ws.Cells.Item(5, 8).Value = "WAB"
$ws.Cells.Item(5, 9).Value = "SEA"
$ws.Cells.Item(5, 10).Value = 1
$ws.Cells.Item(5, 11).Value = 44989
$ws.Cells.Item(5, 12).Value = "DD"
$ws.Cells.Item(5, 13).Value = 12101
$ws.Cells.Item(5, 14).Value = 638
$ws.Cells.Item(5, 15).Value = 2347
$ws.Cells.Item(5, 16).Value = 634
$ws.Cells.Item(5, 17).Value = 2345
$ws.Cells.Item(5, 18).Value = 612
$ws.Cells.Item(5, 19).Value = 2350
$ws.Cells.Item(5, 20).Value = 12101
$ws.Cells.Item(5, 21).Value = "Ngee Ann Poly"

# Row 6
$ws.Cells.Item(6, 1).Value = "NextBus3"
$ws.Cells.Item(6, 2).Value = 74
$ws.Cells.Item(6, 3).Value = 11379
$ws.Cells.Item(6, 4).Value = "Buona Vista Ter"
$ws.Cells.Item(6, 5).Value = "SBST"
$ws.Cells.Item(6, 6).Value = 45684.399375
$ws.Cells.Item(6, 7).Value = 11379
$ws.Cells.Item(6, 8).Value = "WAB"
$ws.Cells.Item(6, 9).Value = "SEA"
$ws.Cells.Item(6, 10).Value = 1
$ws.Cells.Item(6, 11).Value = 64009
$ws.Cells.Item(6, 12).Value = "DD"
$ws.Cells.Item(6, 13).Value = 12109
$ws.Cells.Item(6, 14).Value = 605
$ws.Cells.Item(6, 15).Value = 30
$ws.Cells.Item(6, 16).Value = 559
$ws.Cells.Item(6, 17).Value = 24
$ws.Cells.Item(6, 18).Value = 609
$ws.Cells.Item(6, 19).Value = 27
$ws.Cells.Item(6, 20).Value = 12109
$ws.Cells.Item(6, 21).Value = "Opp Ngee Ann Poly"

# Row 7
$ws.Cells.Item(7, 1).Value = "NextBus3"
$ws.Cells.Item(7, 2).Value = 61
$ws.Cells.Item(7, 3).Value = 82009
$ws.Cells.Item(7, 4).Value = "Eunos Int"
$ws.Cells.Item(7, 5).Value = "SMRT"
$ws.Cells.Item(7, 6).Value = 45684.39452546297
$ws.Cells.Item(7, 7).Value = 82009
$ws.Cells.Item(7, 8).Value = "WAB"
$ws.Cells.Item(7, 9).Value = "SEA"
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 11).Value = 43009
$ws.Cells.Item(7, 12).Value = "SD"
$ws.Cells.Item(7, 13).Value = 12109
$ws.Cells.Item(7, 14).Value = 538
$ws.Cells.Item(7, 15).Value = 2349
$ws.Cells.Item(7, 16).Value = 538
$ws.Cells.Item(7, 17).Value = 2349
$ws.Cells.Item(7, 18).Value = 538
$ws.Cells.Item(7, 19).Value = 2349
$ws.Cells.Item(7, 20).Value = 12109
$ws.Cells.Item(7, 21).Value = "Opp Ngee Ann Poly"

# Row 8
$ws.Cells.Item(8, 1).Value = "NextBus3"
$ws.Cells.Item(8, 2).Value = 154
$ws.Cells.Item(8, 3).Value = 82009
$ws.Cells.Item(8, 4).Value = "Eunos Int"
$ws.Cells.Item(8, 5).Value = "SBST"
$ws.Cells.Item(8, 6).Value = 45684.39274305556
$ws.Cells.Item(8, 7).Value = 82009
$ws.Cells.Item(8, 8).Value = "WAB"
$ws.Cells.Item(8, 9).Value = "SDA"
$ws.Cells.Item(8, 10).Value = 1
$ws.Cells.Item(8, 11).Value = 22009
$ws.Cells.Item(8, 12).Value = "SD"
$ws.Cells.Item(8, 13).Value = 12101
$ws.Cells.Item(8, 14).Value = 601
$ws.Cells.Item(8, 15).Value = 57
$ws.Cells.Item(8, 16).Value = 637
$ws.Cells.Item(8, 17).Value = 55
$ws.Cells.Item(8, 18).Value = 604
$ws.Cells.Item(8, 19).Value = 53
$ws.Cells.Item(8, 20).Value = 12101
$ws.Cells.Item(8, 21).Value = "Ngee Ann Poly"

# Row 9
$ws.Cells.Item(9, 1).Value = "NextBus3"
$ws.Cells.Item(9, 2).Value = 154
$ws.Cells.Item(9, 3).Value = 22009
$ws.Cells.Item(9, 4).Value = "Boon Lay Int"
$ws.Cells.Item(9, 5).Value = "SBST"
$ws.Cells.Item(9, 6).Value = 45684.3990162037
$ws.Cells.Item(9, 7).Value = 22009
$ws.Cells.Item(9, 8).Value = "WAB"
$ws.Cells.Item(9, 9).Value = "SEA"
$ws.Cells.Item(9, 10).Value = 1
$ws.Cells.Item(9, 11).Value = 82009
$ws.Cells.Item(9, 12).Value = "SD"
$ws.Cells.Item(9, 13).Value = 12109
$ws.Cells.Item(9, 14).Value = 546
$ws.Cells.Item(9, 15).Value = 16
$ws.Cells.Item(9, 16).Value = 616
$ws.Cells.Item(9, 17).Value = 15
$ws.Cells.Item(9, 18).Value = 547
$ws.Cells.Item(9, 19).Value = 15
$ws.Cells.Item(9, 20).Value = 12109
$ws.Cells.Item(9, 21).Value = "Opp Ngee Ann Poly"

# Row 10
$ws.Cells.Item(10, 1).Value = "NextBus3"
$ws.Cells.Item(10, 2).Value = 61
$ws.Cells.Item(10, 3).Value = 43009
$ws.Cells.Item(10, 4).Value = "Bt Batok Int"
$ws.Cells.Item(10, 5).Value = "SMRT"
$ws.Cells.Item(10, 6).Value = 45684.40625
$ws.Cells.Item(10, 7).Value = 43009
$ws.Cells.Item(10, 8).Value = "WAB"
$ws.Cells.Item(10, 9).Value = "SDA"
$ws.Cells.Item(10, 10).Value = 1
$ws.Cells.Item(10, 11).Value = 82009
$ws.Cells.Item(10, 12).Value = "SD"
$ws.Cells.Item(10, 13).Value = 12101
$ws.Cells.Item(10, 14).Value = 645
$ws.Cells.Item(10, 15).Value = 108
$ws.Cells.Item(10, 16).Value = 652
$ws.Cells.Item(10, 17).Value = 110
$ws.Cells.Item(10, 18).Value = 642
$ws.Cells.Item(10, 19).Value = 109
$ws.Cells.Item(10, 20).Value = 12101
$ws.Cells.Item(10, 21).Value = "Ngee Ann Poly"

# Row 11
$ws.Cells.Item(11, 1).Value = "NextBus3"
$ws.Cells.Item(11, 2).Value = 74
$ws.Cells.Item(11, 3).Value = 64009
$ws.Cells.Item(11, 4).Value = "Hougang Ctrl Int"
$ws.Cells.Item(11, 5).Value = "SBST"
$ws.Cells.Item(11, 6).Value = 45684.39164351852
$ws.Cells.Item(11, 7).Value = 64009
$ws.Cells.Item(11, 8).Value = "WAB"
$ws.Cells.Item(11, 9).Value = "SEA"
$ws.Cells.Item(11, 10).Value = 1
$ws.Cells.Item(11, 11).Value = 11379
$ws.Cells.Item(11, 12).Value = "SD"
$ws.Cells.Item(11, 13).Value = 12101
$ws.Cells.Item(11, 14).Value = 602
$ws.Cells.Item(11, 15).Value = 2343
$ws.Cells.Item(11, 16).Value = 623
$ws.Cells.Item(11, 17).Value = 2340
$ws.Cells.Item(11, 18).Value = 549
$ws.Cells.Item(11, 19).Value = 2343
$ws.Cells.Item(11, 20).Value = 12101
$ws.Cells.Item(11, 21).Value = "Ngee Ann Poly"

# Row 12
$ws.Cells.Item(12, 1).Value = "NextBus3"
$ws.Cells.Item(12, 2).Value = 151
$ws.Cells.Item(12, 3).Value = 64009
$ws.Cells.Item(12, 4).Value = "Hougang Ctrl Int"
$ws.Cells.Item(12, 5).Value = "SBST"
$ws.Cells.Item(12, 6).Value = 45684.40273148148
$ws.Cells.Item(12, 7).Value = 64009
$ws.Cells.Item(12, 8).Value = "WAB"
$ws.Cells.Item(12, 9).Value = "SEA"
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 11).Value = 16009
$ws.Cells.Item(12, 12).Value = "SD"
$ws.Cells.Item(12, 13).Value = 12101
$ws.Cells.Item(12, 14).Value = 642
$ws.Cells.Item(12, 15).Value = 2317
$ws.Cells.Item(12, 16).Value = 642
$ws.Cells.Item(12, 17).Value = 2315
$ws.Cells.Item(12, 18).Value = 611
$ws.Cells.Item(12, 19).Value = 2317
$ws.Cells.Item(12, 20).Value = 12101
$ws.Cells.Item(12, 21).Value = "Ngee Ann Poly"

# Row 13
$ws.Cells.Item(13, 1).Value = "NextBus3"
$ws.Cells.Item(13, 2).Value = 52
$ws.Cells.Item(13, 3).Value = 28009
$ws.Cells.Item(13, 4).Value = "Jurong East Int"
$ws.Cells.Item(13, 5).Value = "SBST"
$ws.Cells.Item(13, 6).Value = 45684.39863425926
$ws.Cells.Item(13, 7).Value = 28009
$ws.Cells.Item(13, 8).Value = "WAB"
$ws.Cells.Item(13, 9).Value = "SEA"
$ws.Cells.Item(13, 10).Value = 1
$ws.Cells.Item(13, 11).Value = 53009
$ws.Cells.Item(13, 12).Value = "SD"
$ws.Cells.Item(13, 13).Value = 12109
$ws.Cells.Item(13, 14).Value = 623
$ws.Cells.Item(13, 15).Value = 23
$ws.Cells.Item(13, 16).Value = 625
$ws.Cells.Item(13, 17).Value = 21
$ws.Cells.Item(13, 18).Value = 627
$ws.Cells.Item(13, 19).Value = 22
$ws.Cells.Item(13, 20).Value = 12109
$ws.Cells.Item(13, 21).Value = "Opp Ngee Ann Poly"

# Row 14
$ws.Cells.Item(14, 1).Value = "NextBus3"
$ws.Cells.Item(14, 2).Value = 151
$ws.Cells.Item(14, 3).Value = 16009
$ws.Cells.Item(14, 4).Value = "Kent Ridge Ter"
$ws.Cells.Item(14, 5).Value = "SBST"
$ws.Cells.Item(14, 6).Value = 45684.38590277778
$ws.Cells.Item(14, 7).Value = 16009
$ws.Cells.Item(14, 8).Value = "WAB"
$ws.Cells.Item(14, 9).Value = "SEA"
$ws.Cells.Item(14, 10).Value = 1
$ws.Cells.Item(14, 11).Value = 64009
$ws.Cells.Item(14, 12).Value = "SD"
$ws.Cells.Item(14, 13).Value = 12109
$ws.Cells.Item(14, 14).Value = 635
$ws.Cells.Item(14, 15).Value = 2347
$ws.Cells.Item(14, 16).Value = 634
$ws.Cells.Item(14, 17).Value = 2351
$ws.Cells.Item(14, 18).Value = 639
$ws.Cells.Item(14, 19).Value = 2354
$ws.Cells.Item(14, 20).Value = 12109
$ws.Cells.Item(14, 21).Value = "Opp Ngee Ann Poly"

# Row 15
$ws.Cells.Item(15, 1).Value = "NextBus3"
$ws.Cells.Item(15, 2).Value = 75
$ws.Cells.Item(15, 3).Value = 10009
$ws.Cells.Item(15, 4).Value = "Bt Merah Int"
$ws.Cells.Item(15, 5).Value = "SMRT"
$ws.Cells.Item(15, 6).Value = 45684.3955324074
$ws.Cells.Item(15, 7).Value = 10009
$ws.Cells.Item(15, 8).Value = "WAB"
$ws.Cells.Item(15, 9).Value = "SEA"
$ws.Cells.Item(15, 10).Value = 1
$ws.Cells.Item(15, 11).Value = 44989
$ws.Cells.Item(15, 12).Value = "SD"
$ws.Cells.Item(15, 13).Value = 12109
$ws.Cells.Item(15, 14).Value = 548
$ws.Cells.Item(15, 15).Value = 2350
$ws.Cells.Item(15, 16).Value = 546
$ws.Cells.Item(15, 17).Value = 2350
$ws.Cells.Item(15, 18).Value = 552
$ws.Cells.Item(15, 19).Value = 2351
$ws.Cells.Item(15, 20).Value = 12109
$ws.Cells.Item(15, 21).Value = "Opp Ngee Ann Poly"

# Preserve/apply datetime number format for column F (EstimatedTimeOfArrival)
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
